$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F27").Value = 76
$ws.Range("G27").Value = 2725.36
$ws.Range("B34").Value = 63511.97
$ws.Range("F36").Value = 109
$ws.Range("G36").Value = 21447.93
$ws.Range("F38").Value = 467
$ws.Range("G38").Value = 17003.47
$ws.Range("F41").Value = 233
$ws.Range("G41").Value = 44943.37
$ws.Range("F48").Value = 248
$ws.Range("G48").Value = 13912.8
$ws.Range("F51").Value = 164
$ws.Range("G51").Value = 15340.56
$ws.Range("F55").Value = 136
$ws.Range("G55").Value = 7583.36
$ws.Range("F58").Value = 87
$ws.Range("G58").Value = 6779.91
$ws.Range("F61").Value = 253
$ws.Range("G61").Value = 65964.69
$ws.Range("B66").Value = 229553.98
$ws.Range("F72").Value = 11
$ws.Range("G72").Value = 2530.88
$ws.Range("B83").Value = 20340.84
$ws.Range("F112").Value = 42
$ws.Range("G112").Value = 3317.16
$ws.Range("F113").Value = 112
$ws.Range("G113").Value = 15719.2
$ws.Range("B123").Value = 78225.95
$ws.Range("F142").Value = 16
$ws.Range("G142").Value = 1588.96
$ws.Range("B147").Value = 23911.47
$ws.Range("F154").Value = 315
$ws.Range("G154").Value = 10502.1
$ws.Range("B155").Value = 39998.39
$ws.Range("B161").Value = 64350
$ws.Range("E161").Value = 70.63
$ws.Range("F161").Value = 2
$ws.Range("G161").Value = 132.88
$ws.Range("B162").Value = 57756
$ws.Range("E162").Value = 79.37
$ws.Range("F162").Value = -100
$ws.Range("G162").Value = -6644
$ws.Range("F172").Value = 112
$ws.Range("G172").Value = 7114.24
$ws.Range("F173").Value = 58
$ws.Range("G173").Value = 4559.38
$ws.Range("F177").Value = 281
$ws.Range("G177").Value = 13086.17
$ws.Range("F179").Value = 4
$ws.Range("G179").Value = 136.44
$ws.Range("F180").Value = 45
$ws.Range("G180").Value = 7579.8
$ws.Range("F182").Value = 33
$ws.Range("G182").Value = 2954.82
$ws.Range("F186").Value = 43
$ws.Range("G186").Value = 1861.04
$ws.Range("F192").Value = 28
$ws.Range("G192").Value = 1676.92
$ws.Range("B193").Value = 73488.31
$ws.Range("F212").Value = 81
$ws.Range("G212").Value = 7216.29
$ws.Range("F216").Value = 96
$ws.Range("G216").Value = 7132.8
$ws.Range("B218").Value = 88436.08
$ws.Range("F222").Value = 1310
$ws.Range("G222").Value = 24235
$ws.Range("F223").Value = 48
$ws.Range("G223").Value = 1028.64
$ws.Range("F227").Value = 60
$ws.Range("G227").Value = 6876
$ws.Range("B229").Value = 38831.48
$ws.Range("F261").Value = 85
$ws.Range("G261").Value = 6970.85
$ws.Range("F265").Value = 52
$ws.Range("G265").Value = 4635.8
$ws.Range("F267").Value = 148
$ws.Range("G267").Value = 6287.04
$ws.Range("F269").Value = 13
$ws.Range("G269").Value = 1114.36
$ws.Range("F274").Value = 61
$ws.Range("G274").Value = 2125.24
$ws.Range("F277").Value = 15
$ws.Range("G277").Value = 756.15
$ws.Range("F278").Value = 53
$ws.Range("G278").Value = 7183.62
$ws.Range("F280").Value = 25
$ws.Range("G280").Value = 2427.5
$ws.Range("F285").Value = 15
$ws.Range("G285").Value = 1664.1
$ws.Range("F287").Value = 75
$ws.Range("G287").Value = 4105.5
$ws.Range("B295").Value = 139935.99
$ws.Range("F325").Value = 63
$ws.Range("G325").Value = 9523.71
$ws.Range("B328").Value = 7894.81
$ws.Range("F349").Value = 154
$ws.Range("G349").Value = 11493.02
$ws.Range("F351").Value = 222
$ws.Range("G351").Value = 32110.08
$ws.Range("B356").Value = 81131.07
$ws.Range("F361").Value = 275
$ws.Range("G361").Value = 38662.25
$ws.Range("F362").Value = 29
$ws.Range("G362").Value = 21579.19
$ws.Range("B363").Value = 85120.92
$ws.Range("F365").Value = 20
$ws.Range("G365").Value = 1106.6
$ws.Range("F370").Value = 258
$ws.Range("G370").Value = 42825.42
$ws.Range("B372").Value = 70923.53
$ws.Range("F387").Value = 474
$ws.Range("G387").Value = 45788.4
$ws.Range("B389").Value = 63024.43
$ws.Range("F402").Value = 69
$ws.Range("G402").Value = 2367.39
$ws.Range("F408").Value = 233
$ws.Range("G408").Value = 3693.05
$ws.Range("F413").Value = 98
$ws.Range("G413").Value = 5641.86
$ws.Range("F414").Value = 71
$ws.Range("G414").Value = 3542.19
$ws.Range("B417").Value = 184842.57
$ws.Range("F423").Value = 55
$ws.Range("G423").Value = 4597.45
$ws.Range("B427").Value = 24145.39
$ws.Range("F453").Value = 42
$ws.Range("G453").Value = 6103.44
$ws.Range("B458").Value = 110034.49
$ws.Range("F503").Value = 6
$ws.Range("G503").Value = 591
$ws.Range("B508").Value = 41915.65
$ws.Range("F523").Value = 168
$ws.Range("G523").Value = 14382.48
$ws.Range("F524").Value = 28
$ws.Range("G524").Value = 2473.52
$ws.Range("B525").Value = 137330.81
$ws.Range("F528").Value = 316
$ws.Range("G528").Value = 5011.76
$ws.Range("F529").Value = 133
$ws.Range("G529").Value = 4403.63
$ws.Range("F530").Value = 42
$ws.Range("G530").Value = 1813.56
$ws.Range("B535").Value = 29177.77
$ws.Range("F544").Value = 53
$ws.Range("G544").Value = 3280.7
$ws.Range("B556").Value = 69098.68
$ws.Range("F558").Value = 247
$ws.Range("G558").Value = 30096.95
$ws.Range("B561").Value = 35816.41
$ws.Range("F568").Value = 3
$ws.Range("G568").Value = 753.3
$ws.Range("F569").Value = 6
$ws.Range("G569").Value = 3508.32
$ws.Range("B573").Value = 34502.97
$ws.Range("F609").Value = 43
$ws.Range("G609").Value = 4678.83
$ws.Range("F615").Value = 114
$ws.Range("G615").Value = 17632.38
$ws.Range("F616").Value = 15
$ws.Range("G616").Value = 2141.55
$ws.Range("F617").Value = 46
$ws.Range("G617").Value = 2213.52
$ws.Range("F620").Value = 386
$ws.Range("G620").Value = 30335.74
$ws.Range("F622").Value = 503
$ws.Range("G622").Value = 51763.73
$ws.Range("F625").Value = 355
$ws.Range("G625").Value = 13074.65
$ws.Range("F626").Value = 19
$ws.Range("G626").Value = 896.99
$ws.Range("B628").Value = 234038.54
$ws.Range("F659").Value = 53
$ws.Range("G659").Value = 2837.62
$ws.Range("F662").Value = 57
$ws.Range("G662").Value = 4577.67
$ws.Range("F664").Value = 114
$ws.Range("G664").Value = 1026
$ws.Range("B668").Value = 14940.55
$ws.Range("F674").Value = 1070
$ws.Range("G674").Value = 174527.7
$ws.Range("B680").Value = 179802.03
$ws.Range("F687").Value = 2
$ws.Range("G687").Value = 81.98
$ws.Range("B691").Value = 11914.35
$ws.Range("B718").Value = 3157390.5
$ws.Range("B719").Value = 3157390.5
